$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 49 currently carries the special "last row" formatting (a filled-style
# variant on columns B/C). That distinctive formatting needs to move down to
# the new last row (50), so copy it there first...
$ws.Range("B49:C49").Copy()
$ws.Range("B50:C50").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)      # xlPasteFormats

# ...then normalize row 49 back to the regular (non-special) row style used
# by the rest of the visible entries above it (e.g. row 48).
$ws.Range("B48:C48").Copy()
$ws.Range("B49:C49").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# New work-log entry.
$ws.Range("A50").Value = 45756
$ws.Range("B50").Value = 4
$ws.Range("C50").Value = "Finalized report structure and format, finalized all the codes"

Write-Host "Added row 50 work-log entry"
